# Backlog and other things updated for the start of the sprint
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark several backlog rows as done ("x") in column D
$ws.Range("D5").Value = "x"
$ws.Range("D8").Value = "x"
$ws.Range("D10").Value = "x"
$ws.Range("D11").Value = "x"
$ws.Range("D17").Value = "x"

# Row 9 gets upgraded from a single "x" to "xx"
$ws.Range("D9").Value = "xx"

# Leave the selection where the user left off at the start of the sprint
$ws.Range("A24").Select()
